$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Rows($ws, $row1, $row2, $firstCol, $lastCol) {
    $range1 = $ws.Range("$firstCol$row1`:$lastCol$row1")
    $range2 = $ws.Range("$firstCol$row2`:$lastCol$row2")
    $temp = $range1.Value2
    $range1.Value2 = $range2.Value2
    $range2.Value2 = $temp
}

# Swap the data (columns B through AC) between row 11 and row 12
Swap-Rows $ws 11 12 "B" "AC"

# Swap the data (columns B through AC) between row 83 and row 84
Swap-Rows $ws 83 84 "B" "AC"
